# Updates currentAveragePrice / profit figures across the Leve-profit sheets
# (market-price refresh from the scheduled data-pull runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 16433.334
$ws.Range("J51").Value = 16433.334
$ws.Range("L51").Value = 16433.334
$ws.Range("N51").Value = -17401.334

$ws.Range("H64").Value = 14284.071
$ws.Range("I64").Value = 9998.286
$ws.Range("J64").Value = 18569.857
$ws.Range("K64").Value = 9998.286
$ws.Range("L64").Value = 18569.857
$ws.Range("M64").Value = -9750.286
$ws.Range("N64").Value = -19065.857

$ws.Range("H67").Value = 14284.071
$ws.Range("I67").Value = 9998.286
$ws.Range("J67").Value = 18569.857
$ws.Range("K67").Value = 9998.286
$ws.Range("L67").Value = 18569.857
$ws.Range("M67").Value = -9140.286
$ws.Range("N67").Value = -20285.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 5266.6665
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 7700
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 7700
$ws.Range("M11").Value = -256
$ws.Range("N11").Value = -7988

$ws.Range("H32").Value = 5919.65
$ws.Range("I32").Value = 5061.6313
$ws.Range("K32").Value = 5061.6313
$ws.Range("M32").Value = -4774.6313

$ws.Range("H63").Value = 4649.9
$ws.Range("I63").Value = 1499.909
$ws.Range("K63").Value = 1499.909
$ws.Range("M63").Value = -813.9090000000001

$ws.Range("H66").Value = 4649.9
$ws.Range("I66").Value = 1499.909
$ws.Range("K66").Value = 7499.545
$ws.Range("M66").Value = -4067.545

$ws.Range("H74").Value = 1817.5
$ws.Range("I74").Value = 1801.375
$ws.Range("J74").Value = 1882
$ws.Range("K74").Value = 1801.375
$ws.Range("L74").Value = 1882
$ws.Range("M74").Value = -927.375
$ws.Range("N74").Value = -3630

$ws.Range("H77").Value = 1817.5
$ws.Range("I77").Value = 1801.375
$ws.Range("J77").Value = 1882
$ws.Range("K77").Value = 9006.875
$ws.Range("L77").Value = 9410
$ws.Range("M77").Value = -4638.875
$ws.Range("N77").Value = -18146

$ws.Range("H88").Value = 5809.9
$ws.Range("I88").Value = 4251
$ws.Range("K88").Value = 4251
$ws.Range("M88").Value = -3845

$ws.Range("H91").Value = 5809.9
$ws.Range("I91").Value = 4251
$ws.Range("K91").Value = 4251
$ws.Range("M91").Value = -2847

$ws.Range("H122").Value = 2503.5
$ws.Range("I122").Value = 2503.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7510.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5060.5
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5409.3
$ws.Range("I105").Value = 4148.25
$ws.Range("K105").Value = 4148.25
$ws.Range("M105").Value = -2401.25

$ws.Range("H134").Value = 4997
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1299.1666
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650

$ws.Range("H93").Value = 22499.334
$ws.Range("I93").Value = 22499.334
$ws.Range("K93").Value = 22499.334
$ws.Range("M93").Value = -20627.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 143
$ws.Range("I11").Value = 143
$ws.Range("K11").Value = 429
$ws.Range("M11").Value = -289

$ws.Range("H51").Value = 1165.75
$ws.Range("I51").Value = 1165.75
$ws.Range("K51").Value = 3497.25
$ws.Range("M51").Value = -3037.25

$ws.Range("H55").Value = 1901
$ws.Range("I55").Value = 402.5
$ws.Range("K55").Value = 1207.5
$ws.Range("M55").Value = -1030.5

$ws.Range("H107").Value = 1547.5834
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 1732.1
$ws.Range("K107").Value = 1875
$ws.Range("L107").Value = 5196.299999999999
$ws.Range("M107").Value = 45
$ws.Range("N107").Value = -9036.299999999999

$ws.Range("H131").Value = 1557
$ws.Range("I131").Value = 1000
$ws.Range("K131").Value = 3000
$ws.Range("M131").Value = 2040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("N17").Value = -2336

$ws.Range("H55").Value = 9998
$ws.Range("I55").Value = 9998
$ws.Range("K55").Value = 9998
$ws.Range("M55").Value = -9671

$ws.Range("H70").Value = 7555.3
$ws.Range("I70").Value = 8919.727999999999
$ws.Range("K70").Value = 8919.727999999999
$ws.Range("M70").Value = -8649.727999999999

$ws.Range("H73").Value = 7555.3
$ws.Range("I73").Value = 8919.727999999999
$ws.Range("K73").Value = 8919.727999999999
$ws.Range("M73").Value = -7983.727999999999

$ws.Range("H122").Value = 1997
$ws.Range("I122").Value = 1997
$ws.Range("K122").Value = 5991
$ws.Range("M122").Value = -3541

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1917.7
$ws.Range("I22").Value = 1034.1666
$ws.Range("K22").Value = 1034.1666
$ws.Range("M22").Value = -739.1666

$ws.Range("H27").Value = 1917.7
$ws.Range("I27").Value = 1034.1666
$ws.Range("K27").Value = 1034.1666
$ws.Range("M27").Value = -927.1666

$ws.Range("H46").Value = 3153.4
$ws.Range("I46").Value = 2450
$ws.Range("J46").Value = 3409.182
$ws.Range("K46").Value = 2450
$ws.Range("L46").Value = 3409.182
$ws.Range("M46").Value = -2262
$ws.Range("N46").Value = -3785.182

$ws.Range("H55").Value = 660.63635
$ws.Range("I55").Value = 270.22223
$ws.Range("K55").Value = 270.22223
$ws.Range("M55").Value = -97.22223000000002

$ws.Range("H132").Value = 4425.0835
$ws.Range("I132").Value = 4022.3333
$ws.Range("K132").Value = 12066.9999
$ws.Range("M132").Value = -9536.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2256.3572
$ws.Range("I132").Value = 2313.2222
$ws.Range("K132").Value = 6939.6666
$ws.Range("M132").Value = -4409.6666

$ws.Range("H136").Value = 3668
$ws.Range("I136").Value = 3323
$ws.Range("K136").Value = 9969
$ws.Range("M136").Value = -7419
